$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 44, shifting existing rows 44-150 down to 45-151
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the values from the commit diff.
# Most columns (A,B,C,E,F,G,H,I,N,O,Q,R) mirror the row template already
# present (copied down by Insert), only D,J,K,L,M,P differ per the diff.
$ws.Range("A44").Value = 3
$ws.Range("B44").Value = "Femacal de La Calera"
$ws.Range("C44").Value = "Coquimbo"
$ws.Range("D44").Value = 44622
$ws.Range("E44").Value = 5
$ws.Range("F44").Value = 100112052
$ws.Range("G44").Value = "Albahaca"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 110
$ws.Range("K44").Value = 4500
$ws.Range("L44").Value = 5000
$ws.Range("M44").Value = 4727
$ws.Range("N44").Value = "$/docena de matas"
$ws.Range("O44").Value = "Provincia de Quillota"
$ws.Range("P44").Value = 788
$ws.Range("Q44").Value = 6
$ws.Range("R44").Value = "Hortaliza"
